# Final edits - code cleanup
# Rewrites the three role/concept tables: "Roles" columns become
# "Stereotypes" columns, role values are reshuffled (farmer/redneck swap,
# teacher -> laborer/scientists reshuffle, runner/homeless swap), and the
# per-country "F" helper columns are dropped from every sheet (sheet3 keeps
# its styled-but-empty F11:F16 tail).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # White_Concept
$ws2 = $wb.Worksheets.Item(2)   # Asian_Concept
$ws3 = $wb.Worksheets.Item(3)   # Black_Concept

# ---- Header row: Roles -> Stereotypes (same on all three sheets) ----
foreach ($ws in @($ws1, $ws2, $ws3)) {
    $ws.Range("A1").Value = "White_Stereotypes"
    $ws.Range("B1").Value = "Asian_Stereotypes"
    $ws.Range("C1").Value = "Black_Stereotypes"
}

# ---- Role columns A/B/C reshuffle (same pattern on all three sheets) ----
foreach ($ws in @($ws1, $ws2, $ws3)) {
    $ws.Range("A4").Value = "farmer"
    $ws.Range("A7").Value = "redneck"

    $ws.Range("B4").Value = "scientists"
    $ws.Range("B5").Value = "laborer"

    $ws.Range("C4").Value = "criminal"
    $ws.Range("C5").Value = "homeless"
    $ws.Range("C6").Value = "runner"
}

# ---- Drop the per-country "F" columns ----
# sheet1 / sheet2: no data left anywhere in column F -> delete it outright
# (collapses dimension down to A1:D7 and drops rows 8-10 on sheet1).
$ws1.Columns.Item(6).Delete()
$ws2.Columns.Item(6).Delete()

# sheet3: rows 11-16 still hold styled (but empty) F cells, so only clear
# the F1:F7 country values and leave the rest of the sheet untouched.
$ws3.Range("F1:F7").ClearContents()

# ---- Column width tweaks ----
# ColumnWidth uses Excel's character-width units; the engine stores
# width = ColumnWidth + 5/6 rounded to the nearest 1/6, so we back-solve
# for the COM input that reproduces each target stored width.
$ws1.Columns.Item(1).ColumnWidth = 18.666666666666668
$ws1.Columns.Item(2).ColumnWidth = 20.166666666666668

$ws2.Columns.Item(1).ColumnWidth = 23.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 20.166666666666668
$ws2.Columns.Item(3).ColumnWidth = 18.833333333333332

$ws3.Columns.Item(1).ColumnWidth = 18.5
$ws3.Columns.Item(2).ColumnWidth = 19.666666666666668
$ws3.Columns.Item(3).ColumnWidth = 17.166666666666668

# ---- Sheet view / selection / zoom / active sheet ----
# sheet1 becomes the active/selected tab with a new selection + zoom.
$ws2.Range("D13").Select()   # harmless no-op selects before switching away
$ws1.Activate()
$ws1.Range("D10").Select()
$excel.ActiveWindow.Zoom = 142

$ws2.Range("G3:G16").Select()
$ws2.Activate()
$excel.ActiveWindow.Zoom = 100

$ws3.Range("F3:F11").Select()
$ws3.Activate()
$excel.ActiveWindow.Zoom = 214

# Leave sheet1 as the active tab (matches the target workbook view).
$ws1.Activate()
